# "modified test cases on overdue fix"
#
# - Summary sheet: move selection from B11 to C5
# - Repayment schedule sheet: clear out the (now unused) column O values,
#   move selection from A1:P14 to E11, and it is no longer the active tab
# - Transactions sheet: renumber/adjust a handful of transaction rows,
#   drop the two trailing (near-empty) rows 7-8, move selection to D5,
#   and make it the active tab

$wb = $excel.ActiveWorkbook

# --- Summary -----------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
[void]$wsSummary.Range("C5").Select()

# --- Repayment schedule --------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
[void]$wsSchedule.Range("O2:O14").Clear()
[void]$wsSchedule.Range("E11").Select()

# --- Transactions --------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("A2").Value = 1205
$wsTransactions.Range("A3").Value = 1204
$wsTransactions.Range("E3").Value = 92.97
$wsTransactions.Range("A4").Value = 1203
$wsTransactions.Range("E4").Value = 102.94
$wsTransactions.Range("A5").Value = 1196

# Remove the two trailing placeholder rows (old rows 7 & 8)
[void]$wsTransactions.Rows.Item(7).Delete()
[void]$wsTransactions.Rows.Item(7).Delete()

# Transactions becomes the active sheet / tab, with D5 selected
[void]$wsTransactions.Activate()
[void]$wsTransactions.Range("D5").Select()
